$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "search input" test data
$ws.Name = "search-input-data"

# Remove the second column (header "password" + value "wrongemail123")
# entirely -- the sheet now only needs a single column of data.
$ws.Columns.Item(2).Delete()

# Replace the remaining column's values with the new search-bar test data.
$ws.Range("A1").Value = "searchInput"
$ws.Range("A2").Value = "MacBook"
